$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.465.95"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "3.520.02"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("D7").Value = "3.518.58"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  +3.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.425"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "4.115.12"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "3.519.60"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "67.447.27"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "444.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.624"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000130"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.27%  "
$ws.Range("D26").Value = "3.659.99"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.52%  "
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("E33").Value = "  +4.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").Value = "3.512.70"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "177.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0885"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("E44").Value = "  -3.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.882"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.50%  "
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("E51").Value = "  -2.91%  "
